# "added dalog box for engine record"
# Update the EngineConfigurations sheet with two new measurement columns
# (Red_line_temprature / cycles_per_degree), drop the now-unused EGT column
# from EngineModels, widen the last column on EngineInspectionRecord, and
# leave the EngineInspectionRecord sheet as the active tab/selection.

$wb = $excel.ActiveWorkbook

# --- Sheet "Engines" --------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()
$ws1.Range("C1").Select()

# --- Sheet "EngineModels" ----------------------------------------------
# Remove the obsolete "ExaustGasTemperature" column (C) entirely.
$ws2 = $wb.Worksheets.Item(2)
$ws2.Activate()
$ws2.Range("C1:C5").Clear()
$ws2.Range("C1").Select()

# --- Sheet "EngineConfigurations" --------------------------------------
# Rename the EGT column header and add a new "cycles_per_degree" column,
# and correct the recorded value for row 2.
$ws3 = $wb.Worksheets.Item(3)
$ws3.Activate()
$ws3.Range("E1").Value = "Red_line_temprature"
$ws3.Range("F1").Value = "cycles_per_degree"
$ws3.Range("E2").Value = 858
$ws3.Range("F2").Select()

# --- Sheet "EngineInspectionRecord" -------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Activate()
$ws4.Columns.Item(13).ColumnWidth = 26.35
$ws4.Range("H1").Select()
